$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Min.  Tolong petugas di gate in gate out di training tuk lebih sigap.  Kadang kartu gak bisa buka di gate out,  pas pindah gate sebelah nya baru bisa.  Nah petugas nya malah minta bala bantuan.  Entah ke mana minta nya.'
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B2:J2").Value = $arr
$ws.Range("B2").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Di St.  Hi banyak yg gak bisa gate out.  Silang merah semua.  Kalo saya ke jadian di St.  Senayan dan Blok M,  bisa gate in,  gak bisa gate out.  '
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B3:J3").Value = $arr
$ws.Range("B3").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'AC di stasiun bundaran HI tidak nyala dan panas sekali di dalam peron. Kereta juga mengalami gangguan operasional. Gimana nih MRT Jakarta?Baru operasional udah ada gangguan operasional.'
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B4:J4").Value = $arr
$ws.Range("B4").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'udh telat 10 menit daritadi blom ada yg lewat, yg ke lebak bulus udh ada 4. harusnya 13.26 skrg udh 13.37'
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B5:J5").Value = $arr
$ws.Range("B5").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'ini kali ke-3 saya komplain bau AC MRT yg ga enak seperti bau pesing. 2 komplain saya ga pernah di respon. Padahal yg komplain hal ini juga udh banyak loh. Ada apa sama MRT?'
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B6:J6").Value = $arr
$ws.Range("B6").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Kepada MRT Jakarta, saya sempat ke stasiun MRT bundaran HI tadi malam dan melihat mulai tangga masuk/keluar dan area lantai stasiunnya banyak terkena noda tinta atau tumpahan cat sehingga terlihat kotor dan juga masih ada sisa makanan di lantai eskalator dan kereta.'
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B7:J7").Value = $arr
$ws.Range("B7").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Oh iya mrt bagaimana untuk jawaban kasus saya dimana saya masuk HI saldo e money 150an eh keluar blok A saldo saya 0. Duhh sistem gate nya tolong dong diperbaiki, jangan terus - terus an eror'
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B8:J8").Value = $arr
$ws.Range("B8").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'petugas security kalian tidak pny cukup nyali utk menegur penumpang yg sdh terlihat jelas sdg makan di dlm MRT!!!'
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B9:J9").Value = $arr
$ws.Range("B9").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Ini petugas  @mrtjakarta yang di dalam train kmn ya? Bisa dong menegur penumpang yang pada bawa anak trus seenaknya berisik dan main di atas kursi? Ganggu penumpang yg lain loh.'
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B10:J10").Value = $arr
$ws.Range("B10").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Dear  @mrtjakarta mohon maap, ini kenapa proses beli tiketnya lama banget ya ngantrinya 😂 ada kali satu orang 10 menit mah'
$arr[0,1] = 'Twitter'
$arr[0,2] = 'Negatif'
$arr[0,3] = 'Negatif'
$arr[0,4] = 'Negatif'
$arr[0,5] = 'Negatif'
$arr[0,6] = 0
$arr[0,7] = 4
$arr[0,8] = 'Negatif'
$ws.Range("B11:J11").Value = $arr
$ws.Range("B11").Style = "Normal"

$ws.Range("J12").Value = 'Netral'
$ws.Range("J13").Value = 'Netral'
$ws.Range("J14").Value = 'Netral'
$ws.Range("J15").Value = 'Netral'
$ws.Range("J16").Value = 'Netral'
$ws.Range("J17").Value = 'Netral'
$ws.Range("J18").Value = 'Netral'
$ws.Range("J19").Value = 'Netral'
$ws.Range("J20").Value = 'Netral'
$ws.Range("J21").Value = 'Netral'

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'mantap..sekali sy dr Jawa Timur ikut bangga Jakarta punya MRT. Indonesia maju'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B22:J22").Value = $arr
$ws.Range("B22").Interior.Color = 65535

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Kita patut berbangga, dengan pemerintahan yg sekarang, kita bisa punya MRT, LRT, bandara yang megah. Terharu'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B23:J23").Value = $arr
$ws.Range("B23").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Semua pihak terlibat, semua punya andil masing2 atas beroperasinya MRT Jakarta. Sukses terus MRT Jakarta!!!'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B24:J24").Value = $arr
$ws.Range("B24").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Totalitas dan sukses selalu MRT Jakarta! Semoga fase 2 cepat jadi✨'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B25:J25").Value = $arr
$ws.Range("B25").Interior.Color = 65535

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Selamat telah menyelesaikan fase 1 dan semangat dlm melanjutkan pembangunan ke fase selanjutnya,, dan semoga mimpi saya agar bisa ikut berkontribusi dlm pembangunan MRT bs terwujud'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B26:J26").Value = $arr
$ws.Range("B26").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Kerja bagus MRT 👍'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B27:J27").Value = $arr
$ws.Range("B27").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Terima kasih telah membangun MRT untuk Jakarta yang Lebih Baik'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B28:J28").Value = $arr
$ws.Range("B28").Interior.Color = 65535

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Terima kasih MRT Jakarta telah memberikan pelayanan yang terbaik dan tidak kalah dengan dunia internasional, terus bertambah baik sehingga masyarakat dapat menikmati pelayanan yang diberikan.'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B29:J29").Value = $arr
$ws.Range("B29").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Saya bangga dgn Indonesia. Ini adalah awal Indonesia menuju masa kejayaannya'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B30:J30").Value = $arr
$ws.Range("B30").Style = "Normal"

$arr = New-Object "object[,]" 1,9
$arr[0,0] = 'Mantap...menuju jakarta yang lebih baik'
$arr[0,1] = 'Youtube'
$arr[0,2] = 'Positif'
$arr[0,3] = 'Positif'
$arr[0,4] = 'Positif'
$arr[0,5] = 'Positif'
$arr[0,6] = 4
$arr[0,7] = 0
$arr[0,8] = 'Positif'
$ws.Range("B31:J31").Value = $arr
$ws.Range("B31").Interior.Color = 65535

$ws.Columns.Item(2).ColumnWidth = 249.3333333333333

$ws.Range("B22:J31").Select()
$excel.ActiveWindow.ScrollRow = 5
